$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 22
$ws.Range("L22").Value = 8
$ws.Range("N22").Value = 2
$ws.Range("O22").ClearContents()

# Update row 23
$ws.Range("M23").Clear()
$ws.Range("N23").Value = 8
$ws.Range("O23").Value = 4

# Clear row 24 / 25 extra cells
$ws.Range("O24").Clear()
$ws.Range("O25").Clear()

# Update sheet view: scroll so H1 is top-left cell, and selection at N24
$ws.Range("N24").Select()
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
